$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 2500  # H4: 2765.3333 -> 2500
$ws.Cells.Item(4, 9).Value = 2000  # I4: 2648 -> 2000
$ws.Cells.Item(4, 11).Value = 2000  # K4: 2648 -> 2000
$ws.Cells.Item(4, 13).Value = -1886  # M4: -2534 -> -1886
$ws.Cells.Item(64, 8).Value = 3068  # H64: 2835 -> 3068
$ws.Cells.Item(64, 10).Value = 3500  # J64: 3000 -> 3500
$ws.Cells.Item(64, 12).Value = 3500  # L64: 3000 -> 3500
$ws.Cells.Item(64, 14).Value = -3996  # N64: -3496 -> -3996
$ws.Cells.Item(67, 8).Value = 3068  # H67: 2835 -> 3068
$ws.Cells.Item(67, 10).Value = 3500  # J67: 3000 -> 3500
$ws.Cells.Item(67, 12).Value = 3500  # L67: 3000 -> 3500
$ws.Cells.Item(67, 14).Value = -5216  # N67: -4716 -> -5216
$ws.Cells.Item(70, 8).Value = 22266.666  # H70: 26620 -> 22266.666
$ws.Cells.Item(70, 9).Value = 800  # I70: 950 -> 800
$ws.Cells.Item(70, 11).Value = 2400  # K70: 2850 -> 2400
$ws.Cells.Item(70, 13).Value = -2130  # M70: -2580 -> -2130
$ws.Cells.Item(73, 8).Value = 22266.666  # H73: 26620 -> 22266.666
$ws.Cells.Item(73, 9).Value = 800  # I73: 950 -> 800
$ws.Cells.Item(73, 11).Value = 2400  # K73: 2850 -> 2400
$ws.Cells.Item(73, 13).Value = -1464  # M73: -1914 -> -1464
$ws.Cells.Item(74, 8).Value = 3380.375  # H74: 3449.2222 -> 3380.375
$ws.Cells.Item(74, 10).Value = 4999.6665  # J74: 4749.75 -> 4999.6665
$ws.Cells.Item(74, 12).Value = 4999.6665  # L74: 4749.75 -> 4999.6665
$ws.Cells.Item(74, 14).Value = -6871.6665  # N74: -6621.75 -> -6871.6665
$ws.Cells.Item(77, 8).Value = 3380.375  # H77: 3449.2222 -> 3380.375
$ws.Cells.Item(77, 10).Value = 4999.6665  # J77: 4749.75 -> 4999.6665
$ws.Cells.Item(77, 12).Value = 24998.3325  # L77: 23748.75 -> 24998.3325
$ws.Cells.Item(77, 14).Value = -34358.3325  # N77: -33108.75 -> -34358.3325
$ws.Cells.Item(98, 8).Value = 1228.4193  # H98: 1304.7667 -> 1228.4193
$ws.Cells.Item(98, 9).Value = 1030.5769  # I98: 1106.5416 -> 1030.5769
$ws.Cells.Item(98, 10).Value = 2257.2  # J98: 2097.6667 -> 2257.2
$ws.Cells.Item(98, 11).Value = 1030.5769  # K98: 1106.5416 -> 1030.5769
$ws.Cells.Item(98, 12).Value = 2257.2  # L98: 2097.6667 -> 2257.2
$ws.Cells.Item(98, 13).Value = 467.4231  # M98: 391.4584 -> 467.4231
$ws.Cells.Item(98, 14).Value = -5253.2  # N98: -5093.6667 -> -5253.2
$ws.Cells.Item(106, 8).Value = 2358.3845  # H106: 2365.6155 -> 2358.3845
$ws.Cells.Item(106, 9).Value = 2560.4546  # I106: 2569 -> 2560.4546
$ws.Cells.Item(106, 11).Value = 2560.4546  # K106: 2569 -> 2560.4546
$ws.Cells.Item(106, 13).Value = -1929.4546  # M106: -1938 -> -1929.4546
$ws.Cells.Item(122, 8).Value = 1228.4193  # H122: 1304.7667 -> 1228.4193
$ws.Cells.Item(122, 9).Value = 1030.5769  # I122: 1106.5416 -> 1030.5769
$ws.Cells.Item(122, 10).Value = 2257.2  # J122: 2097.6667 -> 2257.2
$ws.Cells.Item(122, 11).Value = 3091.7307  # K122: 3319.6248 -> 3091.7307
$ws.Cells.Item(122, 12).Value = 6771.599999999999  # L122: 6293.000100000001 -> 6771.599999999999
$ws.Cells.Item(122, 13).Value = -641.7307000000001  # M122: -869.6248000000001 -> -641.7307000000001
$ws.Cells.Item(122, 14).Value = -11671.6  # N122: -11193.0001 -> -11671.6
$ws.Cells.Item(139, 8).Value = 69991.664  # H139: 69995 -> 69991.664
$ws.Cells.Item(139, 10).Value = 69991.664  # J139: 69995 -> 69991.664
$ws.Cells.Item(139, 12).Value = 69991.664  # L139: 69995 -> 69991.664
$ws.Cells.Item(139, 14).Value = -80271.664  # N139: -80275 -> -80271.664
$ws.Cells.Item(140, 8).Value = 79076.766  # H140: 81036.07000000001 -> 79076.766
$ws.Cells.Item(140, 10).Value = 79076.766  # J140: 81036.07000000001 -> 79076.766
$ws.Cells.Item(140, 12).Value = 79076.766  # L140: 81036.07000000001 -> 79076.766
$ws.Cells.Item(140, 14).Value = -89436.766  # N140: -91396.07000000001 -> -89436.766

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1504.9286  # H45: 1528.7858 -> 1504.9286
$ws.Cells.Item(45, 9).Value = 875.2  # I45: 942 -> 875.2
$ws.Cells.Item(45, 11).Value = 875.2  # K45: 942 -> 875.2
$ws.Cells.Item(45, 13).Value = -498.2  # M45: -565 -> -498.2
$ws.Cells.Item(61, 8).Value = 7583.2104  # H61: 7223.4 -> 7583.2104
$ws.Cells.Item(61, 9).Value = 9424.546  # I61: 8671.416999999999 -> 9424.546
$ws.Cells.Item(61, 11).Value = 9424.546  # K61: 8671.416999999999 -> 9424.546
$ws.Cells.Item(61, 13).Value = -9212.546  # M61: -8459.416999999999 -> -9212.546
$ws.Cells.Item(88, 8).Value = 4367.6665  # H88: 3580.4285 -> 4367.6665
$ws.Cells.Item(88, 9).Value = 1753  # I88: 2600 -> 1753
$ws.Cells.Item(88, 10).Value = 5675  # J88: 3972.6 -> 5675
$ws.Cells.Item(88, 11).Value = 1753  # K88: 2600 -> 1753
$ws.Cells.Item(88, 12).Value = 5675  # L88: 3972.6 -> 5675
$ws.Cells.Item(88, 13).Value = -1347  # M88: -2194 -> -1347
$ws.Cells.Item(88, 14).Value = -6487  # N88: -4784.6 -> -6487
$ws.Cells.Item(91, 8).Value = 4367.6665  # H91: 3580.4285 -> 4367.6665
$ws.Cells.Item(91, 9).Value = 1753  # I91: 2600 -> 1753
$ws.Cells.Item(91, 10).Value = 5675  # J91: 3972.6 -> 5675
$ws.Cells.Item(91, 11).Value = 1753  # K91: 2600 -> 1753
$ws.Cells.Item(91, 12).Value = 5675  # L91: 3972.6 -> 5675
$ws.Cells.Item(91, 13).Value = -349  # M91: -1196 -> -349
$ws.Cells.Item(91, 14).Value = -8483  # N91: -6780.6 -> -8483
$ws.Cells.Item(97, 8).Value = 1189.1428  # H97: 1172.9546 -> 1189.1428
$ws.Cells.Item(97, 9).Value = 1109.7222  # I97: 1095.1578 -> 1109.7222
$ws.Cells.Item(97, 11).Value = 1109.7222  # K97: 1095.1578 -> 1109.7222
$ws.Cells.Item(97, 13).Value = -613.7221999999999  # M97: -599.1578 -> -613.7221999999999
$ws.Cells.Item(122, 8).Value = 1090.1305  # H122: 1167.3334 -> 1090.1305
$ws.Cells.Item(122, 9).Value = 731.4706  # I122: 816.8570999999999 -> 731.4706
$ws.Cells.Item(122, 10).Value = 2106.3333  # J122: 1868.2858 -> 2106.3333
$ws.Cells.Item(122, 11).Value = 2194.4118  # K122: 2450.5713 -> 2194.4118
$ws.Cells.Item(122, 12).Value = 6318.999899999999  # L122: 5604.857400000001 -> 6318.999899999999
$ws.Cells.Item(122, 13).Value = 255.5882000000001  # M122: -0.57129999999961 -> 255.5882000000001
$ws.Cells.Item(122, 14).Value = -11218.9999  # N122: -10504.8574 -> -11218.9999
$ws.Cells.Item(132, 8).Value = 1660.1666  # H132: 1625.1316 -> 1660.1666
$ws.Cells.Item(132, 9).Value = 1355.3572  # I132: 1331.3 -> 1355.3572
$ws.Cells.Item(132, 11).Value = 4066.0716  # K132: 3993.9 -> 4066.0716
$ws.Cells.Item(132, 13).Value = -1536.0716  # M132: -1463.9 -> -1536.0716
$ws.Cells.Item(136, 8).Value = 7583.2104  # H136: 7223.4 -> 7583.2104
$ws.Cells.Item(136, 9).Value = 9424.546  # I136: 8671.416999999999 -> 9424.546
$ws.Cells.Item(136, 11).Value = 28273.638  # K136: 26014.251 -> 28273.638
$ws.Cells.Item(136, 13).Value = -25723.638  # M136: -23464.251 -> -25723.638

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 204655.8  # H86: 128647.31 -> 204655.8
$ws.Cells.Item(86, 9).Value = 5319.75  # I86: 4113.0835 -> 5319.75
$ws.Cells.Item(86, 10).Value = 1002000  # J86: 502250 -> 1002000
$ws.Cells.Item(86, 11).Value = 5319.75  # K86: 4113.0835 -> 5319.75
$ws.Cells.Item(86, 12).Value = 1002000  # L86: 502250 -> 1002000
$ws.Cells.Item(86, 13).Value = -4196.75  # M86: -2990.0835 -> -4196.75
$ws.Cells.Item(86, 14).Value = -1004246  # N86: -504496 -> -1004246
$ws.Cells.Item(89, 8).Value = 204655.8  # H89: 128647.31 -> 204655.8
$ws.Cells.Item(89, 9).Value = 5319.75  # I89: 4113.0835 -> 5319.75
$ws.Cells.Item(89, 10).Value = 1002000  # J89: 502250 -> 1002000
$ws.Cells.Item(89, 11).Value = 26598.75  # K89: 20565.4175 -> 26598.75
$ws.Cells.Item(89, 12).Value = 5010000  # L89: 2511250 -> 5010000
$ws.Cells.Item(89, 13).Value = -20982.75  # M89: -14949.4175 -> -20982.75
$ws.Cells.Item(89, 14).Value = -5021232  # N89: -2522482 -> -5021232
$ws.Cells.Item(107, 8).Value = 947.44446  # H107: 950.1111 -> 947.44446
$ws.Cells.Item(107, 9).Value = 722.8889  # I107: 728.2222 -> 722.8889
$ws.Cells.Item(107, 11).Value = 722.8889  # K107: 728.2222 -> 722.8889
$ws.Cells.Item(107, 13).Value = 1197.1111  # M107: 1191.7778 -> 1197.1111
$ws.Cells.Item(130, 8).Value = 36363.637  # H130: 37000 -> 36363.637
$ws.Cells.Item(134, 8).Value = 4428.7427  # H134: 4428.857 -> 4428.7427
$ws.Cells.Item(134, 9).Value = 4838.2144  # I134: 4838.357 -> 4838.2144
$ws.Cells.Item(134, 11).Value = 14514.6432  # K134: 14515.071 -> 14514.6432
$ws.Cells.Item(134, 13).Value = -11979.6432  # M134: -11980.071 -> -11979.6432

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3128.1428  # H132: 3134.5 -> 3128.1428
$ws.Cells.Item(132, 9).Value = 1848.25  # I132: 1859.375 -> 1848.25
$ws.Cells.Item(132, 11).Value = 5544.75  # K132: 5578.125 -> 5544.75
$ws.Cells.Item(132, 13).Value = -3014.75  # M132: -3048.125 -> -3014.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 11720.466  # H131: 18127.957 -> 11720.466
$ws.Cells.Item(131, 10).Value = 12535.794  # J131: 20261.61 -> 12535.794
$ws.Cells.Item(131, 12).Value = 37607.382  # L131: 60784.83 -> 37607.382
$ws.Cells.Item(131, 14).Value = -47687.382  # N131: -70864.83 -> -47687.382

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 0  # H62: 29110 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 29110 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 29110 -> 0
$ws.Cells.Item(62, 14).ClearContents()  # N62: was -30482
$ws.Cells.Item(65, 8).Value = 0  # H65: 29110 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 29110 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 87330 -> 0
$ws.Cells.Item(65, 14).ClearContents()  # N65: was -94194
$ws.Cells.Item(97, 8).Value = 628.9677  # H97: 615.69696 -> 628.9677
$ws.Cells.Item(97, 9).Value = 630  # I97: 614.8276 -> 630
$ws.Cells.Item(97, 11).Value = 630  # K97: 614.8276 -> 630
$ws.Cells.Item(97, 13).Value = -134  # M97: -118.8276 -> -134
$ws.Cells.Item(113, 8).Value = 1531.1666  # H113: 1495.6 -> 1531.1666
$ws.Cells.Item(113, 10).Value = 1586.6  # J113: 1556 -> 1586.6
$ws.Cells.Item(113, 12).Value = 1586.6  # L113: 1556 -> 1586.6
$ws.Cells.Item(113, 14).Value = -5926.6  # N113: -5896 -> -5926.6
$ws.Cells.Item(122, 8).Value = 2119.6428  # H122: 2148.1428 -> 2119.6428
$ws.Cells.Item(122, 9).Value = 2096.7  # I122: 2207.4443 -> 2096.7
$ws.Cells.Item(122, 10).Value = 2177  # J122: 2041.4 -> 2177
$ws.Cells.Item(122, 11).Value = 6290.099999999999  # K122: 6622.3329 -> 6290.099999999999
$ws.Cells.Item(122, 12).Value = 6531  # L122: 6124.200000000001 -> 6531
$ws.Cells.Item(122, 13).Value = -3840.099999999999  # M122: -4172.3329 -> -3840.099999999999
$ws.Cells.Item(122, 14).Value = -11431  # N122: -11024.2 -> -11431
$ws.Cells.Item(126, 8).Value = 2695802  # H126: 3144611.8 -> 2695802
$ws.Cells.Item(126, 9).Value = 5053432.5  # I126: 5558442.5 -> 5053432.5
$ws.Cells.Item(126, 10).Value = 102408.8  # J126: 127323.625 -> 102408.8
$ws.Cells.Item(126, 11).Value = 15160297.5  # K126: 16675327.5 -> 15160297.5
$ws.Cells.Item(126, 12).Value = 307226.4  # L126: 381970.875 -> 307226.4
$ws.Cells.Item(126, 13).Value = -15157827.5  # M126: -16672857.5 -> -15157827.5
$ws.Cells.Item(126, 14).Value = -312166.4  # N126: -386910.875 -> -312166.4
$ws.Cells.Item(132, 8).Value = 1480713.4  # H132: 1539905.2 -> 1480713.4
$ws.Cells.Item(132, 9).Value = 2565167.8  # I132: 2748327.8 -> 2565167.8
$ws.Cells.Item(132, 10).Value = 1912.091  # J132: 1913 -> 1912.091
$ws.Cells.Item(132, 11).Value = 7695503.399999999  # K132: 8244983.399999999 -> 7695503.399999999
$ws.Cells.Item(132, 12).Value = 5736.272999999999  # L132: 5739 -> 5736.272999999999
$ws.Cells.Item(132, 13).Value = -7692973.399999999  # M132: -8242453.399999999 -> -7692973.399999999
$ws.Cells.Item(132, 14).Value = -10796.273  # N132: -10799 -> -10796.273

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2307.75  # H7: 2396.72 -> 2307.75
$ws.Cells.Item(7, 9).Value = 2109.0417  # I7: 2178.0908 -> 2109.0417
$ws.Cells.Item(7, 10).Value = 3500  # J7: 4000 -> 3500
$ws.Cells.Item(7, 11).Value = 2109.0417  # K7: 2178.0908 -> 2109.0417
$ws.Cells.Item(7, 12).Value = 3500  # L7: 4000 -> 3500
$ws.Cells.Item(7, 13).Value = -1997.0417  # M7: -2066.0908 -> -1997.0417
$ws.Cells.Item(7, 14).Value = -3724  # N7: -4224 -> -3724
$ws.Cells.Item(22, 8).Value = 2229.5  # H22: 2294.4 -> 2229.5
$ws.Cells.Item(22, 9).Value = 1087.25  # I22: 1249.5 -> 1087.25
$ws.Cells.Item(22, 11).Value = 1087.25  # K22: 1249.5 -> 1087.25
$ws.Cells.Item(22, 13).Value = -792.25  # M22: -954.5 -> -792.25
$ws.Cells.Item(27, 8).Value = 2229.5  # H27: 2294.4 -> 2229.5
$ws.Cells.Item(27, 9).Value = 1087.25  # I27: 1249.5 -> 1087.25
$ws.Cells.Item(27, 11).Value = 1087.25  # K27: 1249.5 -> 1087.25
$ws.Cells.Item(27, 13).Value = -980.25  # M27: -1142.5 -> -980.25
$ws.Cells.Item(40, 8).Value = 8045.409  # H40: 8290.956 -> 8045.409
$ws.Cells.Item(40, 9).Value = 7373.9473  # I40: 7384.1055 -> 7373.9473
$ws.Cells.Item(40, 10).Value = 12298  # J40: 12598.5 -> 12298
$ws.Cells.Item(40, 11).Value = 7373.9473  # K40: 7384.1055 -> 7373.9473
$ws.Cells.Item(40, 12).Value = 12298  # L40: 12598.5 -> 12298
$ws.Cells.Item(40, 13).Value = -7237.9473  # M40: -7248.1055 -> -7237.9473
$ws.Cells.Item(40, 14).Value = -12570  # N40: -12870.5 -> -12570
$ws.Cells.Item(46, 8).Value = 1822.4286  # H46: 1818.1428 -> 1822.4286
$ws.Cells.Item(46, 9).Value = 1111.4286  # I46: 1085 -> 1111.4286
$ws.Cells.Item(46, 10).Value = 2533.4285  # J46: 2795.6667 -> 2533.4285
$ws.Cells.Item(46, 11).Value = 1111.4286  # K46: 1085 -> 1111.4286
$ws.Cells.Item(46, 12).Value = 2533.4285  # L46: 2795.6667 -> 2533.4285
$ws.Cells.Item(46, 13).Value = -923.4286  # M46: -897 -> -923.4286
$ws.Cells.Item(46, 14).Value = -2909.4285  # N46: -3171.6667 -> -2909.4285
$ws.Cells.Item(55, 8).Value = 614  # H55: 626.6 -> 614
$ws.Cells.Item(55, 10).Value = 596  # J55: 617.375 -> 596
$ws.Cells.Item(55, 12).Value = 596  # L55: 617.375 -> 596
$ws.Cells.Item(55, 14).Value = -942  # N55: -963.375 -> -942
$ws.Cells.Item(61, 8).Value = 2502.9524  # H61: 2506.762 -> 2502.9524
$ws.Cells.Item(61, 10).Value = 2451.25  # J61: 2461.25 -> 2451.25
$ws.Cells.Item(61, 12).Value = 2451.25  # L61: 2461.25 -> 2451.25
$ws.Cells.Item(61, 14).Value = -2855.25  # N61: -2865.25 -> -2855.25
$ws.Cells.Item(93, 8).Value = 25642146  # H93: 30304018 -> 25642146
$ws.Cells.Item(93, 9).Value = 877.5  # I93: 885.8889 -> 877.5
$ws.Cells.Item(93, 10).Value = 111113040  # J93: 166668110 -> 111113040
$ws.Cells.Item(93, 11).Value = 877.5  # K93: 885.8889 -> 877.5
$ws.Cells.Item(93, 12).Value = 111113040  # L93: 166668110 -> 111113040
$ws.Cells.Item(93, 13).Value = 370.5  # M93: 362.1111 -> 370.5
$ws.Cells.Item(93, 14).Value = -111115536  # N93: -166670606 -> -111115536
$ws.Cells.Item(113, 8).Value = 2502.9524  # H113: 2506.762 -> 2502.9524
$ws.Cells.Item(113, 10).Value = 2451.25  # J113: 2461.25 -> 2451.25
$ws.Cells.Item(113, 12).Value = 2451.25  # L113: 2461.25 -> 2451.25
$ws.Cells.Item(113, 14).Value = -6791.25  # N113: -6801.25 -> -6791.25
$ws.Cells.Item(126, 8).Value = 2307.75  # H126: 2396.72 -> 2307.75
$ws.Cells.Item(126, 9).Value = 2109.0417  # I126: 2178.0908 -> 2109.0417
$ws.Cells.Item(126, 10).Value = 3500  # J126: 4000 -> 3500
$ws.Cells.Item(126, 11).Value = 6327.125100000001  # K126: 6534.2724 -> 6327.125100000001
$ws.Cells.Item(126, 12).Value = 10500  # L126: 12000 -> 10500
$ws.Cells.Item(126, 13).Value = -3857.125100000001  # M126: -4064.2724 -> -3857.125100000001
$ws.Cells.Item(126, 14).Value = -15440  # N126: -16940 -> -15440
$ws.Cells.Item(132, 8).Value = 1717.9546  # H132: 1717.9773 -> 1717.9546
$ws.Cells.Item(132, 9).Value = 1178.2858  # I132: 1178.3214 -> 1178.2858
$ws.Cells.Item(132, 11).Value = 3534.8574  # K132: 3534.9642 -> 3534.8574
$ws.Cells.Item(132, 13).Value = -1004.8574  # M132: -1004.9642 -> -1004.8574

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 697.13336  # H107: 661.6 -> 697.13336
$ws.Cells.Item(107, 9).Value = 521.1667  # I107: 476.75 -> 521.1667
$ws.Cells.Item(107, 11).Value = 1563.5001  # K107: 1430.25 -> 1563.5001
$ws.Cells.Item(107, 13).Value = 356.4999  # M107: 489.75 -> 356.4999
$ws.Cells.Item(126, 8).Value = 1457.2727  # H126: 1645.3158 -> 1457.2727
$ws.Cells.Item(126, 9).Value = 1342.5555  # I126: 1557.8 -> 1342.5555
$ws.Cells.Item(126, 11).Value = 4027.6665  # K126: 4673.4 -> 4027.6665
$ws.Cells.Item(126, 13).Value = -1557.6665  # M126: -2203.4 -> -1557.6665
$ws.Cells.Item(132, 8).Value = 1405.8474  # H132: 1449.2982 -> 1405.8474
$ws.Cells.Item(132, 9).Value = 1021.5476  # I132: 1053.975 -> 1021.5476
$ws.Cells.Item(132, 10).Value = 2355.2942  # J132: 2379.4707 -> 2355.2942
$ws.Cells.Item(132, 11).Value = 3064.6428  # K132: 3161.925 -> 3064.6428
$ws.Cells.Item(132, 12).Value = 7065.882599999999  # L132: 7138.4121 -> 7065.882599999999
$ws.Cells.Item(132, 13).Value = -534.6428000000001  # M132: -631.9249999999997 -> -534.6428000000001
$ws.Cells.Item(132, 14).Value = -12125.8826  # N132: -12198.4121 -> -12125.8826
$ws.Cells.Item(136, 8).Value = 1159  # H136: 1159.8064 -> 1159
$ws.Cells.Item(136, 9).Value = 860.76  # I136: 861.76 -> 860.76
$ws.Cells.Item(136, 11).Value = 2582.28  # K136: 2585.28 -> 2582.28
$ws.Cells.Item(136, 13).Value = -32.27999999999975  # M136: -35.27999999999975 -> -32.27999999999975
